$wb = $excel.ActiveWorkbook

# Nordex is the 2nd sheet (tab index 1, 0-based) and gets a new "MpulseID" header
# in column B next to the existing "WorkOrderID" header in column A.
$ws = $wb.Worksheets.Item("Nordex")
$ws.Range("B1").Value = "MpulseID"

# Make Nordex the active sheet/tab, with B1 selected.
$ws.Activate()
$ws.Range("B1").Select()
